$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the model parameter values (column B, rows 2-13)
$ws.Range("B2").Value = 0.1914962479904777
$ws.Range("B3").Value = -0.1142643096335048
$ws.Range("B4").Value = -0.1475071333334574
$ws.Range("B5").Value = -0.2231475834934862
$ws.Range("B6").Value = -0.4342047056573254
$ws.Range("B7").Value = 0.1232856110648965
$ws.Range("B8").Value = 0.255648357988816
$ws.Range("B9").Value = 0.001991875202858949
$ws.Range("B10").Value = 0.0646952354760759
$ws.Range("B11").Value = 0.8732820434956737
$ws.Range("B12").Value = 0.8384664120137821
$ws.Range("B13").Value = -0.03003781585575057

# Remove the now-obsolete last row (shot_during_regular_play)
$ws.Rows.Item(14).Delete()
